$wb = $excel.ActiveWorkbook

$oldGuid = "e43303ca-72c4-43e4-8774-cc44c9b48287"
$newGuid = "ac687601-052b-44d0-972f-b586b7726b98"

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4ff03ed2b24b92892c6a3dbe912542566285b30a/e2e/$oldGuid.md"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-23 21:01:20"

# The hyperlink object on external links is read-only in-place here, so
# drop it and re-add one that keeps the same target URL (still the old
# guid, matching upstream) but refreshes the displayed text.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $baseUrl, [Type]::Missing, [Type]::Missing, "e2e\$newGuid.md")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.8e70bb869513f033c2b586b83a01bfb1ad5a25f0.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-23 21:01:15"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $baseUrl, [Type]::Missing, [Type]::Missing, "$newGuid.md")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.8e70bb869513f033c2b586b83a01bfb1ad5a25f0.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-23 21:01:20"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $baseUrl, [Type]::Missing, [Type]::Missing, "$newGuid.md")
